$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values on specific rows to match re-pulled data / mean calculation
$ws.Range("F2").Value = -1
$ws.Range("F4").Value = 1
$ws.Range("F5").Value = 1
$ws.Range("F6").Value = -3
$ws.Range("F9").Value = -4
$ws.Range("F13").Value = 3
$ws.Range("F14").Value = -10
$ws.Range("F15").Value = -5
$ws.Range("F19").Value = -1
